$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Add hours spent for sprint 9 of the first week
$ws.Range("I5").Value = 11
$ws.Range("H6").Value = 4
$ws.Range("I7").Value = 4
$ws.Range("F8").Value = 8
$ws.Range("I8").Value = 30
$ws.Range("B9").Value = 8
$ws.Range("E9").Value = 1
$ws.Range("H9").Value = 11
$ws.Range("C10").Value = 12

# Leave the final selection on H9, as in the authored edit
$ws.Range("H9").Select()
